$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Tests_AOS"

# New column L header
$ws.Range("L1").Value = "Test 10"

# Row 14 - AOS Existing Account credentials
$ws.Range("K14").Value = "elad1234"
$ws.Range("L14").Value = "elad1234"

# Row 15 - password
$ws.Range("K15").Value = "Thbyrby145"
$ws.Range("L15").Value = "Thbyrby145"

# Row 16 - AOS New Account username
$ws.Range("J16").Value = "ga12347"

# Row 17 - mail
$ws.Range("J17").Value = "ga12l@gmail.com"

# Row 18 - password
$ws.Range("J18").Value = "ga2341"

# Row 19 - SafePay User username
$ws.Range("J19").Value = "eg12345678"

# Row 20 - password
$ws.Range("J20").Value = "Eg2345"

# Row 21-24 hold purely-numeric looking text, force them to be stored as text
# (matching the source data which is inlineStr, not numeric) by pre-formatting
# the range as Text before assigning the values.
$ws.Range("K21:K24").NumberFormat = "@"

# Row 21 - MasterCredit Card
$ws.Range("K21").Value = "123456789123"

# Row 22 - CVV
$ws.Range("K22").Value = "123"

# Row 23 - Month
$ws.Range("K23").Value = "2"

# Row 24 - Year
$ws.Range("K24").Value = "3"

# Row 25 - Name
$ws.Range("K25").Value = "elad-gal"

# Row 26 - Test Result columns J and L (K26 already has "V")
$ws.Range("J26").Value = "V"
$ws.Range("L26").Value = "V"
